# Generate Report for Handback
# Updates the "zh-cn" and "de-de" status sheets: row 5 (a38377e6-...) now has
# a handback file reported, but it is not the latest handback version, so an
# error message is recorded and the handback datetime/target columns are
# widened to fit the new content.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e714b08542f8096a1a81cd0d807b6dba63bd084d/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09156c08591adc94e50adc97c6316b5606dc5d74/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09156c08591adc94e50adc97c6316b5606dc5d74/e2e/a38377e6-a599-44f8-87da-f903eaf54708.md"

$sheets = @(
    @{ Name = "zh-cn"; TargetFile = "a38377e6-a599-44f8-87da-f903eaf54708.6f9b2d4a208bcb3b5061e774984fb615200f46c9.zh-cn.xlf"; HandbackTime = "2016-08-30 08:16:38" },
    @{ Name = "de-de"; TargetFile = "a38377e6-a599-44f8-87da-f903eaf54708.6f9b2d4a208bcb3b5061e774984fb615200f46c9.de-de.xlf"; HandbackTime = "2016-08-30 08:16:56" }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen the columns that now hold longer content (Latest Target File,
    # Latest Handback File, Error Detail) to match the other 40-wide columns.
    $refWidth = $ws.Columns.Item(7).ColumnWidth
    $ws.Columns.Item(9).ColumnWidth = $refWidth
    $ws.Columns.Item(10).ColumnWidth = $refWidth
    $ws.Columns.Item(16).ColumnWidth = $refWidth

    # Row 5 corresponds to a38377e6-a599-44f8-87da-f903eaf54708.md
    # A handback was produced (Latest Target File / Latest Handback File),
    # recorded with a handback datetime, but it does not match the latest
    # handoff, so the row gets an Error Detail message and the handback
    # file name becomes a hyperlink (like the source file name).
    $ws.Hyperlinks.Add($ws.Range("I5"), $latestUrl, "", "", "a38377e6-a599-44f8-87da-f903eaf54708.md") | Out-Null
    $ws.Range("I5").Font.Underline = 2
    $ws.Range("I5").Font.Color = 15570276

    $ws.Range("J5").Value = $info.TargetFile
    $ws.Range("K5").Value = $info.HandbackTime
    $ws.Range("P5").Value = $errorMessage
}
